$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "J2"  = 2.8
    "L2"  = 3.65
    "S2"  = 1.38
    "T2"  = 2.82
    "U2"  = 1.57
    "V2"  = 2.25
    "W2"  = 8.5
    "X2"  = 11.5
    "Z2"  = 22
    "AA2" = 16.5
    "AB2" = 23
    "AE2" = 11.5
    "AH2" = 11.75
    "AI2" = 19.5
    "AJ2" = 10.75
    "AL2" = 26
    "AM2" = 27
    "AO2" = 11.75
    "AP2" = 18
    "AQ2" = 45
    "AR2" = 70
    "AS2" = 200
    "AT2" = 2.82
    "AX2" = 17.5
    "AY2" = 21
    "AZ2" = 80
    "BA2" = 100
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
